$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data values in B2:K16 (leave headers in row 1 and labels in column A untouched)
$ws.Range("B2:K16").ClearContents()

# Write the new values per the commit's full-series evaluation update
$ws.Range("B2").Value = -0.1340702663479956
$ws.Range("C2").Value = -0.2758992130872696
$ws.Range("D2").Value = 0.3670754083147943
$ws.Range("E2").Value = 0.4914881146443768
$ws.Range("F2").Value = 0.0938579086067595
$ws.Range("G2").Value = 0.01768919490591373
$ws.Range("H2").Value = 0.2940112530288354
$ws.Range("I2").Value = 0.1227914307171113
$ws.Range("J2").Value = 0.2759388362258526
$ws.Range("K2").Value = -0.001489938197266189

$ws.Range("B3").Value = -0.293453334447234
$ws.Range("C3").Value = 0.3495212869548299
$ws.Range("D3").Value = 0.4739339932844123
$ws.Range("E3").Value = 0.07630378724679503
$ws.Range("F3").Value = 0.0001350735459492769
$ws.Range("G3").Value = 0.2764571316688709
$ws.Range("H3").Value = 0.1052373093571469
$ws.Range("I3").Value = 0.2583847148658881
$ws.Range("J3").Value = -0.01904405955723064
$ws.Range("K3").Value = -0.182031752916177

$ws.Range("B4").Value = 0.4376286059058361
$ws.Range("C4").Value = 0.5620413122354185
$ws.Range("D4").Value = 0.1644111061978012
$ws.Range("E4").Value = 0.08824239249695551
$ws.Range("F4").Value = 0.3645644506198771
$ws.Range("G4").Value = 0.1933446283081531
$ws.Range("H4").Value = 0.3464920338168943
$ws.Range("I4").Value = 0.06906325939377558
$ws.Range("J4").Value = -0.0939244339651708
$ws.Range("K4").Value = -0.2180070093596886

$ws.Range("B5").Value = 0.7406159457232021
$ws.Range("C5").Value = 0.3429857396855849
$ws.Range("D5").Value = 0.2668170259847391
$ws.Range("E5").Value = 0.5431390841076607
$ws.Range("F5").Value = 0.3719192617959367
$ws.Range("G5").Value = 0.525066667304678
$ws.Range("H5").Value = 0.2476378928815592
$ws.Range("I5").Value = 0.0846501995226128
$ws.Range("J5").Value = -0.03943237587190501
$ws.Range("K5").Value = 0.4767206611340558

$ws.Range("B6").Value = 1.255012967438235
$ws.Range("C6").Value = 1.178844253737389
$ws.Range("D6").Value = 1.455166311860311
$ws.Range("E6").Value = 1.283946489548587
$ws.Range("F6").Value = 1.437093895057328
$ws.Range("G6").Value = 1.159665120634209
$ws.Range("H6").Value = 0.9966774272752628
$ws.Range("I6").Value = 0.8725948518807449
$ws.Range("J6").Value = 1.388747888886706
$ws.Range("K6").Value = 1.178844253737389

$ws.Range("B7").Value = 0.2348700177716323
$ws.Range("C7").Value = 0.5111920758945538
$ws.Range("D7").Value = 0.3399722535828299
$ws.Range("E7").Value = 0.4931196590915711
$ws.Range("F7").Value = 0.2156908846684524
$ws.Range("G7").Value = 0.05270319130950599
$ws.Range("H7").Value = -0.07137938408501182
$ws.Range("I7").Value = 0.444773652920949
$ws.Range("J7").Value = 0.2348700177716323

$ws.Range("B8").Value = 0.5151599734076631
$ws.Range("C8").Value = 0.343940151095939
$ws.Range("D8").Value = 0.4970875566046802
$ws.Range("E8").Value = 0.2196587821815615
$ws.Range("F8").Value = 0.0566710888226151
$ws.Range("G8").Value = -0.06741148657190271
$ws.Range("H8").Value = 0.4487415504340581
$ws.Range("I8").Value = 0.2388379152847414

$ws.Range("B9").Value = 0.4795802412661804
$ws.Range("C9").Value = 0.6327276467749217
$ws.Range("D9").Value = 0.3552988723518029
$ws.Range("E9").Value = 0.1923111789928565
$ws.Range("F9").Value = 0.06822860359833866
$ws.Range("G9").Value = 0.5843816406042994
$ws.Range("H9").Value = 0.3744780054549828

$ws.Range("B10").Value = 0.3919214649192569
$ws.Range("C10").Value = 0.1144926904961382
$ws.Range("D10").Value = -0.04849500286280822
$ws.Range("E10").Value = -0.172577578257326
$ws.Range("F10").Value = 0.3435754587486348
$ws.Range("G10").Value = 0.1336718235993181

$ws.Range("B11").Value = 0.0691614752440418
$ws.Range("C11").Value = -0.09382621811490459
$ws.Range("D11").Value = -0.2179087935094224
$ws.Range("E11").Value = 0.2982442434965384
$ws.Range("F11").Value = 0.08834060834722172

$ws.Range("B12").Value = -0.1606876400509585
$ws.Range("C12").Value = -0.2847702154454763
$ws.Range("D12").Value = 0.2313828215604846
$ws.Range("E12").Value = 0.02147918641116785

$ws.Range("B13").Value = -0.3143564178021929
$ws.Range("C13").Value = 0.201796619203768
$ws.Range("D13").Value = -0.00810701594554874

$ws.Range("B14").Value = 0.1836459624741271
$ws.Range("C14").Value = -0.02625767267518964

$ws.Range("B15").Value = -0.04428949692388896


